$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) New sheet "Turns per day" appended after the last existing sheet.
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "Turns per day"

# ---------------------------------------------------------------------------
# 2) Populate cells. Order matters: it controls the shared-string table
#    insertion order, so text is written in the exact sequence the strings
#    must first appear to land on the right shared-string indices.
# ---------------------------------------------------------------------------
$ws.Range("A2").Value = "STEPS_PER_REVOLUTION"
$ws.Range("A1").Value = "Param"
$ws.Range("C1").Value = "Unit"
$ws.Range("D1").Value = "Value"
$ws.Range("B1").Value = "Description"
$ws.Range("A3").Value = "WINDER_DELAY_CYCLES"
$ws.Range("C2").Value = "steps"
$ws.Range("H2").Value = "Target TPD"
$ws.Range("A4").Value = "WINDER_ROTATIONS_PER_CYCLE"
$ws.Range("C4").Value = "turns"
$ws.Range("A5").Value = "WIND_BOTH_DIRECTIONS"
$ws.Range("H4").Value = "Suggested settings"
$ws.Range("A6").Value = "ms per minute"
$ws.Range("C3").Value = "ms"
$ws.Range("A7").Value = "intervals per day"
$ws.Range("H5").Value = "Intervals per day"
$ws.Range("A8").Value = "time per interval"
$ws.Range("C8").Value = "seconds"
$ws.Range("A9").Value = "time per rotation"
$ws.Range("C5").Value = "boolean"
$ws.Range("C7").Value = "count"
$ws.Range("A10").Value = "ms per day"
$ws.Range("B7").Value = "once every 15min"
$ws.Range("C6").Value = "ms"
$ws.Range("C9").Value = "seconds"
$ws.Range("C10").Value = "ms"

$ws.Range("D2").Value = 2048
$ws.Range("D3").Value = 900000
$ws.Range("D4").Value = 5
$ws.Range("D5").Formula = "=TRUE"
$ws.Range("D6").Value = 60000
$ws.Range("D7").Formula = "=D10/D3"
$ws.Range("D8").Formula = "=D4*D9"
$ws.Range("D9").Value = 4
$ws.Range("D10").Value = 86400000
$ws.Range("I2").Formula = "=D7*D4*IF(D5,2,1)"

# ---------------------------------------------------------------------------
# 3) Styling: bold header row, bold+larger "Target TPD" callout, number
#    formats for the intervals-per-day and ms-per-day figures.
# ---------------------------------------------------------------------------
$ws.Range("A1:D1").Font.Bold = $true
$ws.Range("H2:I2").Font.Bold = $true
$ws.Range("H2:I2").Font.Size = 14
$ws.Range("D10").NumberFormat = "0.00E+00"
$ws.Range("D7").NumberFormat = "0"

# ---------------------------------------------------------------------------
# 4) Column widths / row height matching the authored layout.
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 18.998697916666668
$ws.Columns.Item(2).ColumnWidth = 18.998697916666668
$ws.Columns.Item(8).ColumnWidth = 10.998697916666666
$ws.Rows.Item(2).RowHeight = 19

# ---------------------------------------------------------------------------
# 5) Selection / activation state.
# ---------------------------------------------------------------------------
$ws.Range("I13").Select()
$ws.Activate()

# Tabelle1 selection moved to B31.
$tabelle1 = $wb.Worksheets.Item(1)
$tabelle1.Activate()
$tabelle1.Range("B31").Select()

# Re-activate the new sheet so it ends up the active tab, matching
# activeTab pointing at "Turns per day".
$ws.Activate()
